# Auto-generated Excel COM-interop script to apply cached-value updates
# to the leve profit columns (H-N) across multiple crafting-discipline sheets.
# These cells hold plain numeric values (no formulas in the workbook), so we
# set them directly via Range.Value on each affected cell.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 456.3846
$ws.Range("I2").Value = 264.75
$ws.Range("K2").Value = 264.75
$ws.Range("M2").Value = -151.75

$ws.Range("H17").Value = 603.4146
$ws.Range("J17").Value = 616
$ws.Range("L17").Value = 1848
$ws.Range("N17").Value = -2184

$ws.Range("H29").Value = 2137.5
$ws.Range("J29").Value = 2750
$ws.Range("L29").Value = 8250
$ws.Range("N29").Value = -8812

$ws.Range("H33").Value = 182
$ws.Range("I33").Value = 108.53333
$ws.Range("K33").Value = 108.53333
$ws.Range("M33").Value = 120.46667

$ws.Range("H38").Value = 850.7778
$ws.Range("I38").Value = 332.125
$ws.Range("J38").Value = 5000
$ws.Range("K38").Value = 996.375
$ws.Range("L38").Value = 15000
$ws.Range("M38").Value = -624.375
$ws.Range("N38").Value = -15744

$ws.Range("H43").Value = 6959445
$ws.Range("I43").Value = 22600.2
$ws.Range("J43").Value = 18520852
$ws.Range("K43").Value = 22600.2
$ws.Range("L43").Value = 18520852
$ws.Range("M43").Value = -22531.2
$ws.Range("N43").Value = -18520990

$ws.Range("H58").Value = 750.625
$ws.Range("I58").Value = 667.5
$ws.Range("K58").Value = 2002.5
$ws.Range("M58").Value = -1852.5

$ws.Range("H86").Value = 4136.591
$ws.Range("I86").Value = 3800.5
$ws.Range("J86").Value = 4416.6665
$ws.Range("K86").Value = 3800.5
$ws.Range("L86").Value = 4416.6665
$ws.Range("M86").Value = -2677.5
$ws.Range("N86").Value = -6662.6665

$ws.Range("H87").Value = 41997
$ws.Range("J87").Value = 41997
$ws.Range("L87").Value = 41997
$ws.Range("N87").Value = -44493

$ws.Range("H89").Value = 4136.591
$ws.Range("I89").Value = 3800.5
$ws.Range("J89").Value = 4416.6665
$ws.Range("K89").Value = 19002.5
$ws.Range("L89").Value = 22083.3325
$ws.Range("M89").Value = -13386.5
$ws.Range("N89").Value = -33315.3325

$ws.Range("H90").Value = 41997
$ws.Range("J90").Value = 41997
$ws.Range("L90").Value = 125991
$ws.Range("N90").Value = -138471

$ws.Range("H98").Value = 3529.1428
$ws.Range("I98").Value = 3709.0417
$ws.Range("K98").Value = 3709.0417
$ws.Range("M98").Value = -2211.0417

$ws.Range("H104").Value = 859.3333
$ws.Range("I104").Value = 859.3333
$ws.Range("K104").Value = 2577.9999
$ws.Range("M104").Value = -830.9998999999998

$ws.Range("H112").Value = 2294.5483
$ws.Range("I112").Value = 757.1429000000001
$ws.Range("J112").Value = 2742.9583
$ws.Range("K112").Value = 2271.4287
$ws.Range("L112").Value = 8228.874899999999
$ws.Range("M112").Value = -1163.4287
$ws.Range("N112").Value = -10444.8749

$ws.Range("H121").Value = 1328.4286
$ws.Range("J121").Value = 1315.2307
$ws.Range("L121").Value = 3945.6921
$ws.Range("N121").Value = -7439.6921

$ws.Range("H122").Value = 3529.1428
$ws.Range("I122").Value = 3709.0417
$ws.Range("K122").Value = 11127.1251
$ws.Range("M122").Value = -8677.125100000001

$ws.Range("H138").Value = 1193.57
$ws.Range("J138").Value = 1718.0878
$ws.Range("L138").Value = 5154.2634
$ws.Range("N138").Value = -15434.2634

$ws.Range("H141").Value = 672.8605
$ws.Range("I141").Value = 545.7
$ws.Range("K141").Value = 1637.1
$ws.Range("M141").Value = 3542.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1954.6875
$ws.Range("I110").Value = 1546.2
$ws.Range("J110").Value = 2635.5
$ws.Range("K110").Value = 1546.2
$ws.Range("L110").Value = 2635.5
$ws.Range("M110").Value = 498.8
$ws.Range("N110").Value = -6725.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2685.9
$ws.Range("I20").Value = 2492.8572
$ws.Range("J20").Value = 3136.3333
$ws.Range("K20").Value = 2492.8572
$ws.Range("L20").Value = 3136.3333
$ws.Range("M20").Value = -2245.8572
$ws.Range("N20").Value = -3630.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 16303.75
$ws.Range("I41").Value = 8050
$ws.Range("K41").Value = 8050
$ws.Range("M41").Value = -7622

$ws.Range("H50").Value = 23257
$ws.Range("J50").Value = 23257
$ws.Range("L50").Value = 23257
$ws.Range("N50").Value = -24507

$ws.Range("H51").Value = 17000
$ws.Range("J51").Value = 23000
$ws.Range("L51").Value = 23000
$ws.Range("N51").Value = -24472

$ws.Range("H59").Value = 30000
$ws.Range("J59").Value = 30000
$ws.Range("L59").Value = 30000
$ws.Range("N59").Value = -32290

$ws.Range("H60").Value = 9099.866
$ws.Range("J60").Value = 25499.5
$ws.Range("L60").Value = 25499.5
$ws.Range("N60").Value = -26521.5

$ws.Range("H61").Value = 17000
$ws.Range("J61").Value = 23000
$ws.Range("L61").Value = 23000
$ws.Range("N61").Value = -23696

$ws.Range("H74").Value = 33000
$ws.Range("J74").Value = 33000
$ws.Range("L74").Value = 33000
$ws.Range("N74").Value = -34748

$ws.Range("H77").Value = 33000
$ws.Range("J77").Value = 33000
$ws.Range("L77").Value = 99000
$ws.Range("N77").Value = -107736

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2488.4666
$ws.Range("I81").Value = 1502.4
$ws.Range("K81").Value = 4507.200000000001
$ws.Range("M81").Value = -3384.200000000001

$ws.Range("H84").Value = 2488.4666
$ws.Range("I84").Value = 1502.4
$ws.Range("K84").Value = 13521.6
$ws.Range("M84").Value = -7905.6

$ws.Range("H104").Value = 3786.3157
$ws.Range("J104").Value = 3774
$ws.Range("L104").Value = 11322
$ws.Range("N104").Value = -16564

$ws.Range("H131").Value = 29413212
$ws.Range("I131").Value = 125000350
$ws.Range("J131").Value = 1783.5769
$ws.Range("K131").Value = 375001050
$ws.Range("L131").Value = 5350.7307
$ws.Range("M131").Value = -374996010
$ws.Range("N131").Value = -15430.7307

$ws.Range("H132").Value = 1782
$ws.Range("I132").Value = 1834.8
$ws.Range("J132").Value = 1650
$ws.Range("K132").Value = 16513.2
$ws.Range("L132").Value = 14850
$ws.Range("M132").Value = -13983.2
$ws.Range("N132").Value = -19910

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 151.5
$ws.Range("I2").Value = 171
$ws.Range("J2").Value = 119
$ws.Range("K2").Value = 171
$ws.Range("L2").Value = 119
$ws.Range("M2").Value = -58
$ws.Range("N2").Value = -345

$ws.Range("H63").Value = 15000.111
$ws.Range("J63").Value = 15000.111
$ws.Range("L63").Value = 15000.111
$ws.Range("N63").Value = -16372.111

$ws.Range("H66").Value = 15000.111
$ws.Range("J66").Value = 15000.111
$ws.Range("L66").Value = 45000.333
$ws.Range("N66").Value = -51864.333

$ws.Range("H70").Value = 112502080
$ws.Range("I70").Value = 125002056
$ws.Range("K70").Value = 125002056
$ws.Range("M70").Value = -125001786

$ws.Range("H73").Value = 112502080
$ws.Range("I73").Value = 125002056
$ws.Range("K73").Value = 125002056
$ws.Range("M73").Value = -125001120

$ws.Range("H80").Value = 3599.75
$ws.Range("I80").Value = 2533
$ws.Range("J80").Value = 6800
$ws.Range("K80").Value = 2533
$ws.Range("L80").Value = 6800
$ws.Range("M80").Value = -1535
$ws.Range("N80").Value = -8796

$ws.Range("H83").Value = 3599.75
$ws.Range("I83").Value = 2533
$ws.Range("J83").Value = 6800
$ws.Range("K83").Value = 12665
$ws.Range("L83").Value = 34000
$ws.Range("M83").Value = -7673
$ws.Range("N83").Value = -43984

$ws.Range("H113").Value = 1458.5
$ws.Range("I113").Value = 1604.4
$ws.Range("K113").Value = 1604.4
$ws.Range("M113").Value = 565.5999999999999

$ws.Range("H130").Value = 36493.332
$ws.Range("J130").Value = 36493.332
$ws.Range("L130").Value = 36493.332
$ws.Range("N130").Value = -46533.332

$ws.Range("H133").Value = 50599
$ws.Range("J133").Value = 50599
$ws.Range("L133").Value = 50599
$ws.Range("N133").Value = -60719

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7666.6665
$ws.Range("J46").Value = 7666.6665
$ws.Range("L46").Value = 7666.6665
$ws.Range("N46").Value = -8042.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3177.1292
$ws.Range("I132").Value = 3215.3076
$ws.Range("K132").Value = 9645.9228
$ws.Range("M132").Value = -7115.9228
